$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E contain numeric-looking / percentage-looking text that Excel
# would otherwise auto-convert to numbers. Force them to remain as text by
# setting the NumberFormat to "@" (Text) before assigning the value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '328.16'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.25%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '44.10'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.92%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.571'

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.20%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08050'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.56%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.973'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '4.55%'

$ws.Range("B7").Value = 'MXToken'

$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9459'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.84%'

$ws.Range("B8").Value = 'BTSEToken'

$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.551'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-8.51%'

$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1171'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.00%'

$ws.Range("B10").Value = 'WazirX'

$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1855'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.90%'

$ws.Range("B11").Value = 'MCDex'

$ws.Range("C11").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '11.82'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '38.08%'

$ws.Range("B12").Value = 'MandalaExchangeToken'

$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09787'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.90%'

$ws.Range("B13").Value = 'BitrueCoin'

$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04723'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '13.74%'

$ws.Range("B14").Value = 'BitMartToken'

$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1065'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.45%'

$ws.Range("B15").Value = 'BitForexToken'

$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001284'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.06%'

$ws.Range("B16").Value = 'CoinExToken'

$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04215'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.42%'

$ws.Range("B17").Value = 'TigerCash'

$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.005943'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.57%'

$ws.Range("B18").Value = 'HotbitToken'

$ws.Range("C18").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.004331'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.33%'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.370'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-5.43%'

$ws.Range("B20").Value = 'GateToken'

$ws.Range("C20").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.325'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.94%'

$ws.Range("B21").Value = 'BitpandaEcosystemToken'

$ws.Range("C21").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.3475'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.34%'

$ws.Range("B22").Value = 'ProBitToken'

$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1418'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '4.22%'

$ws.Range("B23").Value = 'ZBToken'

$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.2511'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-3.00%'

$ws.Range("B24").Value = 'BitKan'

$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001252'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.04%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001192'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-3.19%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003750'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-6.06%'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02591'

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-2.58%'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05510'

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '0.63%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007545'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.78%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1400'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.21%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007577'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-33.75%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002019'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.30%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008374'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-9.51%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00007093'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.26%'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.08%'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '1.31%'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '35.97%'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.08%'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.08%'
